# Actualización al 11 de junio de 2023
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Ingreso" (sheet1): add contribution rows for 2023-06-11 (serial 45088)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Ingreso")

$ingresoRows = @(
    @{ Row = 444; Miembro = "Randy";   Aporte = 100 },
    @{ Row = 445; Miembro = "Carlos";  Aporte = 100 },
    @{ Row = 446; Miembro = "Jeicol";  Aporte = 100 },
    @{ Row = 447; Miembro = "Anuel";   Aporte = 100 },
    @{ Row = 448; Miembro = "Joel";    Aporte = 50  },
    @{ Row = 449; Miembro = "Gustavo"; Aporte = 50  },
    @{ Row = 450; Miembro = "Omaury";  Aporte = 100 },
    @{ Row = 451; Miembro = "Kibelo";  Aporte = 800 },
    @{ Row = 452; Miembro = "Robert";  Aporte = 250 },
    @{ Row = 453; Miembro = "Orlando"; Aporte = 500 },
    @{ Row = 454; Miembro = "Michy";   Aporte = 200 }
)

foreach ($item in $ingresoRows) {
    $r = $item.Row
    $prev = $r - 1
    $ws1.Range("A$($prev):D$($prev)").Copy() | Out-Null
    $ws1.Range("A$($r):D$($r)").PasteSpecial(-4122) | Out-Null
    $ws1.Range("A$r").Value2 = 45088
    $ws1.Range("B$r").Value2 = $item.Miembro
    $ws1.Range("C$r").Value2 = $item.Aporte
    $ws1.Range("D$r").Value2 = "Aporte"
}

$excel.CutCopyMode = 0
$ws1.Select()
$ws1.Range("C453").Select()

# ---------------------------------------------------------------------------
# Sheet "Gastos" (sheet2): add two expense rows
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gastos")

$ws2.Range("A44:C44").Copy() | Out-Null
$ws2.Range("A45:C45").PasteSpecial(-4122) | Out-Null
$ws2.Range("A45").Value2 = 45088
$ws2.Range("B45").Value2 = "Agua y hielo"
$ws2.Range("C45").Value2 = 150

$ws2.Range("A44:C44").Copy() | Out-Null
$ws2.Range("A46:C46").PasteSpecial(-4122) | Out-Null
$ws2.Range("A46").Value2 = 45089
$ws2.Range("B46").Value2 = "Neverita"
$ws2.Range("C46").Value2 = 4000

$excel.CutCopyMode = 0
$ws2.Select()
$ws2.Range("A47").Select()

$ws1.Select()
